$d = $word.ActiveDocument

# The document's single section has two distinct headers (primary/default
# and first-page) and two distinct footers (primary/default and
# first-page). Each contains one inline picture whose wp:docPr / pic:cNvPr
# "name" metadata needs to be renamed (Pearson logo: image1.png ->
# image2.png ; BTEC logo: image2.jpg -> image1.jpg). The picture's visible
# content/description is untouched - only the internal file-name label.
#
# InlineShape has no writable "Name" property in the Word object model, so
# each picture is promoted to a floating Shape (which does expose .Name),
# renamed, and then converted back to an inline picture so the layout
# stays exactly as it was (wp:inline, not wp:anchor).

function Rename-InlineShape($range, [string]$newName) {
    $shape = $range.InlineShapes.Item(1).ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

$sec = $d.Sections.Item(1)

# Headers: BTEC logo, image2.jpg -> image1.jpg
Rename-InlineShape $sec.Headers.Item(1) "image1.jpg"   # default header -> header2.xml
Rename-InlineShape $sec.Headers.Item(2) "image1.jpg"   # first-page header -> header1.xml

# Footers: Pearson logo, image1.png -> image2.png
Rename-InlineShape $sec.Footers.Item(1) "image2.png"   # default footer -> footer2.xml
Rename-InlineShape $sec.Footers.Item(2) "image2.png"   # first-page footer -> footer1.xml

Write-Output "Renamed header/footer logo pictures"
